# Reorder / refresh the "Estado de Cuenta" detail rows (17-30).
# The underlying (worker, period) -> (Valor Mora, Salario Basico) pairs are
# unchanged; only the grouping/order changes: rows are now grouped by
# worker (JOSE DAVID ALMEIDA LEONES, then RUBEN DARIO GOMEZ OROZCO) with
# each worker's periods sorted descending (2311 .. 2305).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$tipoDoc = "CC"

$rows = @(
    @{ Row = 17; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2311"; Valor = 112000; Salario = 1160000 },
    @{ Row = 18; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2310"; Valor = 120000; Salario = 1160000 },
    @{ Row = 19; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2309"; Valor = 120000; Salario = 1160000 },
    @{ Row = 20; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2308"; Valor = 120000; Salario = 1160000 },
    @{ Row = 21; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2307"; Valor = 46400;  Salario = 1160000 },
    @{ Row = 22; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2306"; Valor = 46400;  Salario = 1160000 },
    @{ Row = 23; Doc = "73006956";   Nombre = "JOSE DAVID ALMEIDA LEONES"; Periodo = "2305"; Valor = 46400;  Salario = 1160000 },
    @{ Row = 24; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2311"; Valor = 74667;  Salario = 2000000 },
    @{ Row = 25; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2310"; Valor = 80000;  Salario = 2000000 },
    @{ Row = 26; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2309"; Valor = 80000;  Salario = 2000000 },
    @{ Row = 27; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2308"; Valor = 80000;  Salario = 2000000 },
    @{ Row = 28; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2307"; Valor = 80000;  Salario = 2000000 },
    @{ Row = 29; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2306"; Valor = 80000;  Salario = 2000000 },
    @{ Row = 30; Doc = "1143351433"; Nombre = "RUBEN DARIO GOMEZ OROZCO";  Periodo = "2305"; Valor = 80000;  Salario = 2000000 }
)

foreach ($item in $rows) {
    $r = $item.Row
    $ws.Cells.Item($r, 2).Value = $tipoDoc
    $ws.Cells.Item($r, 3).Value = $item.Doc
    $ws.Cells.Item($r, 4).Value = $item.Nombre
    $ws.Cells.Item($r, 5).Value = $item.Periodo
    $ws.Cells.Item($r, 6).Value = $item.Valor
    $ws.Cells.Item($r, 7).Value = $item.Salario
}
